# Generate Report for Handback
# Refresh the handback status for the "15757890-fe2d-41c0-a634-369d4eb07159.md"
# file (row 2 on each per-locale sheet): new handoff/handback timestamps for
# the zh-cn and de-de xliff round-trips, and the rolled-up "Latest HO Xliff
# Generate Date" on the Overview sheet (max across locales).

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-19 14:54:26"
$zhcn.Range("K2").Value = "2016-08-19 14:54:42"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-19 14:54:30"
$dede.Range("K2").Value = "2016-08-19 14:54:49"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-19 14:54:30"
